$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The new day (2020-06-07, serial 43989) repeats the same 20-hospital
# block/layout as the previous day (2020-06-06, rows 1273:1292), just with
# updated occupancy figures. Copy that block down (carries over all cell
# styles/formatting) then patch date + the handful of changed values.
$src = $ws.Range("A1273:H1292")
$dst = $ws.Range("A1293:H1312")
$src.Copy($dst)

# Update the date column for the new block to the next day.
$ws.Range("A1293:A1312").Value = 43989

# Patch the occupancy values that changed vs. the prior day.
$ws.Range("C1293").Value = 7
$ws.Range("C1294").Value = 26
$ws.Range("C1298").Value = 3
$ws.Range("C1300").Value = 8
$ws.Range("C1304").Value = 3
$ws.Range("C1308").ClearContents()
$ws.Range("C1312").ClearContents()
